{"js": "// Append a new \"Hotfix\" section to the end of the status log, following the\n// same pattern as the existing \"Update: ...\" sections in this document: a\n// blank separator paragraph, a header line, then bullet points describing\n// the fix.\n\nconst body = context.document.body;\nlet anchor = body.paragraphs.getLast();\n\nconst newLines = [\n  \"\",\n  \"Hotfix: 2026-02-21 (Team invite CORS + SMTP handling)\",\n  \"- Normalized `AllowedOrigins` parsing by trimming trailing `/` to prevent origin mismatch on Railway/Render URLs.\",\n  \"- Team invite/resend now handle SMTP failures gracefully with explicit 502 response and audit log events instead of unhandled 500.\"\n];\n\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n  // Match the surrounding paragraphs' run formatting (Helvetica Light, 12pt).\n  anchor.font.set({ name: \"Helvetica Light\", size: 12 });\n  if (line === \"\") {\n    // Keep an explicit (empty) text run on the blank separator paragraph,\n    // matching the style of the other blank-line paragraphs already present\n    // in this document.\n    anchor.getRange().insertText(\"\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Append a new \"Hotfix\" section to the end of the status log, following the\n# same pattern as the existing \"Update: ...\" sections in this document: a\n# blank separator paragraph, a header line, then bullet points describing\n# the fix.\n\n$d = $word.ActiveDocument\n\n$lines = @(\n    \"\",\n    \"Hotfix: 2026-02-21 (Team invite CORS + SMTP handling)\",\n    \"- Normalized ``AllowedOrigins`` parsing by trimming trailing ``/`` to prevent origin mismatch on Railway/Render URLs.\",\n    \"- Team invite/resend now handle SMTP failures gracefully with explicit 502 response and audit log events instead of unhandled 500.\"\n)\n\n$r = $d.Paragraphs.Last.Range\n\nforeach ($line in $lines) {\n    $r.Collapse(0)\n    $r.InsertParagraphAfter()\n    $r.Collapse(0)\n    $newPara = $d.Paragraphs.Last\n    $newPara.Range.InsertAfter($line)\n    # Apply font formatting to the text only (not the paragraph mark), so we\n    # don't pollute <w:pPr> with an extra <w:rPr>.\n    $textOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)\n    $textOnly.Font.NameAscii = \"Helvetica Light\"\n    $textOnly.Font.Size = 12\n    $r = $newPara.Range\n}\n"}
